$wb = $excel.ActiveWorkbook

# --- "Repayment schedule" sheet: insert a new (blank) column before the
# existing "Late" column (column N), shifting Late / heading-spacer /
# Outstanding one column to the right (O, P, Q). ---
$ws = $wb.Worksheets.Item("Repayment schedule")
$ws.Columns.Item(14).Insert()

# Give the freshly inserted column N its own explicit width (custom, not
# auto best-fit) instead of inheriting the old "Late" column's best-fit width.
$ws.Columns.Item(14).ColumnWidth = 10.1

# --- Make "Repayment schedule" the active sheet/tab, with S10 selected ---
$ws.Activate()
[void]$ws.Range("S10").Select()
